$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.220.57"
Set-TextValue "E2" "  +0.05%  "
Set-TextValue "D3" "1.901.16"
Set-TextValue "E3" "  +0.38%  "
Set-TextValue "E4" "  -0.12%  "
Set-TextValue "D5" "306.20"
Set-TextValue "E5" "  -0.43%  "
Set-TextValue "E6" "  -0.09%  "
Set-TextValue "D7" "0.5380"
Set-TextValue "E7" "  +3.48%  "
Set-TextValue "E8" "  +1.12%  "
Set-TextValue "D9" "0.07278"
Set-TextValue "E9" "  +0.09%  "
Set-TextValue "D10" "22.15"
Set-TextValue "E10" "  +4.63%  "
Set-TextValue "D11" "0.9027"
Set-TextValue "E11" "  +0.33%  "
Set-TextValue "D12" "0.08194"
Set-TextValue "E12" "  +0.08%  "
Set-TextValue "D13" "96.10"
Set-TextValue "E13" "  -0.47%  "
Set-TextValue "D14" "5.334"
Set-TextValue "E14" "  +0.92%  "
Set-TextValue "E15" "  -0.16%  "
Set-TextValue "E16" "  +2.02%  "
Set-TextValue "D17" "0.000008646"
Set-TextValue "E17" "  +0.56%  "
Set-TextValue "E18" "  -0.09%  "
Set-TextValue "D19" "27.252.37"
Set-TextValue "E19" "  +0.09%  "
Set-TextValue "D20" "5.035"
Set-TextValue "E20" "  -0.98%  "
Set-TextValue "D21" "1.103.48"
Set-TextValue "E21" "  -41.80%  "
Set-TextValue "D22" "10.77"
Set-TextValue "E22" "  +0.83%  "
Set-TextValue "D23" "6.486"
Set-TextValue "E23" "  +1.53%  "
Set-TextValue "D24" "149.60"
Set-TextValue "E24" "  +1.54%  "
Set-TextValue "D25" "2.285"
Set-TextValue "E25" "  -0.91%  "
Set-TextValue "E26" "  +0.71%  "
Set-TextValue "D27" "1.743"
Set-TextValue "E27" "  -0.06%  "
Set-TextValue "D28" "116.55"
Set-TextValue "E28" "  +1.21%  "
Set-TextValue "D29" "4.807"
Set-TextValue "E29" "  -0.24%  "
Set-TextValue "D30" "4.742"
Set-TextValue "E30" "  -3.29%  "
Set-TextValue "D31" "0.09214"
Set-TextValue "E31" "  -0.12%  "
Set-TextValue "D32" "0.8339"
Set-TextValue "E32" "  +4.87%  "
Set-TextValue "D33" "0.05071"
Set-TextValue "E33" "  +0.89%  "
Set-TextValue "D34" "1.215"
Set-TextValue "E34" "  -0.46%  "
Set-TextValue "D35" "2.996"
Set-TextValue "E35" "  +1.37%  "
Set-TextValue "D36" "3.338"
Set-TextValue "E36" "  -3.14%  "
Set-TextValue "D37" "2.688"
Set-TextValue "E37" "  +3.99%  "
Set-TextValue "D38" "0.5806"
Set-TextValue "E38" "  +2.20%  "
Set-TextValue "D39" "0.02002"
Set-TextValue "E39" "  +0.90%  "
Set-TextValue "D40" "1.076"
Set-TextValue "E40" "  +0.12%  "
Set-TextValue "D41" "9.331"
Set-TextValue "E41" "  +4.36%  "
Set-TextValue "D42" "6.597"
Set-TextValue "E42" "  +0.69%  "
Set-TextValue "D43" "116.86"
Set-TextValue "E43" "  +1.24%  "
Set-TextValue "D44" "0.1521"
Set-TextValue "E44" "  +0.56%  "
Set-TextValue "D45" "0.4982"
Set-TextValue "E45" "  +2.49%  "
Set-TextValue "E46" "  -0.13%  "
Set-TextValue "D47" "10.11"
Set-TextValue "E47" "  +0.10%  "
Set-TextValue "E48" "  +1.07%  "
Set-TextValue "D49" "38.43"
Set-TextValue "E49" "  +0.82%  "
Set-TextValue "E50" "  +3.77%  "
Set-TextValue "D51" "63.31"
Set-TextValue "E51" "  -0.14%  "
